# Generate Report for Handback
#
# This mirrors a "handback" localization report refresh:
#   1. Status text "Ready for handoff" -> "Handed back: in sync with en-US"
#      (appears on Overview!E2:F2/E3:F3 and on the zh-cn/de-de Status column C2/C3)
#   2. On the zh-cn and de-de detail sheets, the "Latest Target File" (I) and
#      "Latest Handback File" (J) columns get populated for both data rows,
#      with I2/I3 turned into hyperlinks (like column A already is).
#   3. On de-de, "Latest Handback DateTime" (K2/K3) gets a real timestamp,
#      and the de-de row-3 "Latest Handoff File" (G3) is corrected to the
#      de-de xlf name (it incorrectly held the zh-cn name before).
#   4. Columns C (zh-cn/de-de sheets), I and J get widened to fit the new
#      longer file-name / hyperlink text.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1. Status text refresh (shared string replaced in place everywhere it
#    is used: Overview E2:F2 / E3:F3, zh-cn C2/C3, de-de C2/C3)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. zh-cn sheet: fill Latest Target File (I) / Latest Handback File (J)
#    and refresh the shared "Latest Handback DateTime" placeholder text
#    (the same shared string backs zh-cn K2/K3, so it updates too).
# ---------------------------------------------------------------------
$wsZh.Range("J2").Value = "650f1eed-6b0a-4d34-90d3-c9aa34d8ce0f.21aee6a0276f75c704baf83a82b18715732a4c37.zh-cn.xlf"
$wsZh.Range("J3").Value = "a12e0c7a-e379-4e63-8710-ca2436491d04.965d357403ecdfd151cb12c3f59cccd7b60896c9.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-10-13 13:54:30"
$wsZh.Range("K3").Value = "2016-10-13 13:54:30"

# ---------------------------------------------------------------------
# 3. de-de sheet: fix Latest Handoff File (G3), fill Latest Target File (I),
#    Latest Handback File (J) and Latest Handback DateTime (K)
# ---------------------------------------------------------------------
$wsDe.Range("G3").Value = "a12e0c7a-e379-4e63-8710-ca2436491d04.965d357403ecdfd151cb12c3f59cccd7b60896c9.de-de.xlf"

$wsDe.Range("J2").Value = "650f1eed-6b0a-4d34-90d3-c9aa34d8ce0f.21aee6a0276f75c704baf83a82b18715732a4c37.de-de.xlf"
$wsDe.Range("J3").Value = "a12e0c7a-e379-4e63-8710-ca2436491d04.965d357403ecdfd151cb12c3f59cccd7b60896c9.de-de.xlf"

$wsDe.Range("K2").Value = "2016-10-13 13:54:47"
$wsDe.Range("K3").Value = "2016-10-13 13:54:47"

# ---------------------------------------------------------------------
# 4. Add "Latest Target File" hyperlinks on column I for both sheets,
#    rebuilding the hyperlink collection so the rIds land in the same
#    order as the rest of the workbook (A2, I2, A3, I3). Hyperlinks.Add
#    re-styles its range with the built-in "Hyperlink" theme style, so
#    re-apply the workbook's own custom HyperLink look (underline +
#    FF6495ED) to all four link cells afterwards to match column A.
# ---------------------------------------------------------------------
foreach ($ws in @($wsZh, $wsDe)) {
    $linkA2 = $ws.Hyperlinks.Item(1).Address
    $linkA3 = $ws.Hyperlinks.Item(2).Address

    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $linkA2, "", "", "650f1eed-6b0a-4d34-90d3-c9aa34d8ce0f.md")
    $ws.Range("I2").Value = "650f1eed-6b0a-4d34-90d3-c9aa34d8ce0f.md"
    $ws.Hyperlinks.Add($ws.Range("I2"), $linkA2, "", "", "650f1eed-6b0a-4d34-90d3-c9aa34d8ce0f.md")

    $ws.Hyperlinks.Add($ws.Range("A3"), $linkA3, "", "", "a12e0c7a-e379-4e63-8710-ca2436491d04.md")
    $ws.Range("I3").Value = "a12e0c7a-e379-4e63-8710-ca2436491d04.md"
    $ws.Hyperlinks.Add($ws.Range("I3"), $linkA3, "", "", "a12e0c7a-e379-4e63-8710-ca2436491d04.md")

    foreach ($addr in @("A2", "A3", "I2", "I3")) {
        $r = $ws.Range($addr)
        $r.Font.Underline = 2
        $r.Font.Color = 15570276
    }
}

# ---------------------------------------------------------------------
# 5. Column widths: widen C (17.216 -> 29.978) and I/J (-> 40) on both
#    detail sheets so the longer hyperlink / file-name text is visible.
# ---------------------------------------------------------------------
foreach ($ws in @($wsZh, $wsDe)) {
    $ws.Columns.Item(3).ColumnWidth = 29.9777050018311
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40
}

Write-Output "Generate Report for Handback: done"
